$d = $word.ActiveDocument

# Locate the paragraph that introduces "Asciicam" / "Ascii-image." -- in the
# original document this sentence is split across several runs with
# w:proofErr spell-check bookends around "Asciicam" and "Ascii". We collapse
# it down to a single run with the full sentence and drop the proofErr
# markers (the spelling-check artifacts Word leaves behind are no longer
# needed once the text is merged).
$find = $d.Content
$find.Find.Execute("Asciicam") | Out-Null
$para = $find.Paragraphs(1)
$r = $para.Range

# Build a minimal WordprocessingML package containing just the replacement
# paragraph, preserving the paragraph/run formatting that was already there
# (rsid attributes, en-US language run properties), but as one contiguous
# run -- matching what Word itself does when it folds adjacent same-format
# runs back together and clears the now-stale proofErr bookmarks.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="006877A0" w:rsidRPr="006877A0" w:rsidRDefault="006877A0" w:rsidP="006877A0">
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:r w:rsidRPr="006877A0">
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>Asciicam allows you to convert pictures that you take, or pictures that already exist on your phone into an Ascii-image.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
